$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 475, pushing existing rows 475..504 down to 476..505
$ws.Rows.Item(475).Insert()

# Populate the newly inserted row with the new weekly price record
$ws.Range("A475").Value2 = 3
$ws.Range("B475").Value2 = "Femacal de La Calera"
$ws.Range("C475").Value2 = "Coquimbo"
$ws.Range("D475").Value2 = 44783
$ws.Range("E475").Value2 = 5
$ws.Range("F475").Value2 = 100112003
$ws.Range("G475").Value2 = "Ajo"
$ws.Range("H475").Value2 = "Chino"
$ws.Range("I475").Value2 = "Primera"
$ws.Range("J475").Value2 = 78
$ws.Range("K475").Value2 = 24000
$ws.Range("L475").Value2 = 25000
$ws.Range("M475").Value2 = 24513
$ws.Range("N475").Value2 = "$/caja 10 kilos"
$ws.Range("O475").Value2 = "China"
$ws.Range("P475").Value2 = 2451
$ws.Range("Q475").Value2 = 10
$ws.Range("R475").Value2 = "Hortaliza"
